$d = $word.ActiveDocument

# 1. Merge the split "Output Results to CSV(Avg Fitness + time)" runs
#    (and drop the proofErr markers) into a single plain run via Find/Replace.
$d.Content.Find.Execute(
    "Output Results to CSV(Avg Fitness + time)", $true, $false, $false, $false, $false,
    $true, 1, $false, "Output Results to CSV(Avg Fitness + time)", 2) | Out-Null

# 2. Add "Internal Commentary" text to the paragraph that holds the
#    _GoBack bookmark (inserted right before the bookmark, same paragraph).
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPara.Range.InsertBefore("Internal Commentary")

# 3. Append a new, empty paragraph after that one (before the sectPr).
#    Range.InsertParagraphAfter() always materializes an empty <w:r/> in
#    the freshly-created paragraph, so split that new paragraph again
#    (collapsed-range text assignment does NOT leave a stray run behind)
#    and drop the now-redundant trailing paragraph that carries the artifact.
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkPara.Range.Collapse(0) | Out-Null
$bookmarkPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$splitPoint = $newPara.Range.Start
$d.Range($splitPoint, $splitPoint).Text = [char]13

$artifactPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$artifactPara.Range.Delete() | Out-Null
